$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 33   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/2/2026  Through  2/8/2026"

# --- Simple numeric value updates ---
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -14.285714285714
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 18
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 163.636363636364
$ws.Range("L16").Value = -27.5
$ws.Range("M16").Value = -43.13725490196
$ws.Range("N16").Value = -81.761006289308
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -3.225806451612
$ws.Range("I17").Value = 42
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = 10.526315789473
$ws.Range("L17").Value = -2.325581395348
$ws.Range("M17").Value = 55.555555555555
$ws.Range("N17").Value = -51.162790697674
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 52.941176470588
$ws.Range("L18").Value = -25.714285714285
$ws.Range("M18").Value = -31.578947368421
$ws.Range("N18").Value = -82.191780821917
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 54.054054054054
$ws.Range("I19").Value = 78
$ws.Range("J19").Value = 52
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = 39.285714285714
$ws.Range("M19").Value = 129.411764705882
$ws.Range("N19").Value = 41.818181818181
$ws.Range("C20").Value = 3
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = 66.666666666666
$ws.Range("L20").Value = -31.818181818181
$ws.Range("M20").Value = 7.142857142857
$ws.Range("N20").Value = -86.607142857142
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = 42.424242424242
$ws.Range("I21").Value = 196
$ws.Range("J21").Value = 130
$ws.Range("K21").Value = 50.76923076923
$ws.Range("L21").Value = -2
$ws.Range("M21").Value = 17.365269461077
$ws.Range("N21").Value = -65.794066317626
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 3
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -50
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -33.333333333333
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -11.764705882352
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = 15
$ws.Range("I24").Value = 115
$ws.Range("J24").Value = 110
$ws.Range("K24").Value = 4.545454545454
$ws.Range("L24").Value = 19.791666666666
$ws.Range("M24").Value = 69.117647058823
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 13.333333333333
$ws.Range("I25").Value = 22
$ws.Range("J25").Value = 23
$ws.Range("K25").Value = -4.347826086956
$ws.Range("L25").Value = 37.5
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = -2
$ws.Range("I26").Value = 66
$ws.Range("J26").Value = 71
$ws.Range("K26").Value = -7.042253521126
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 6.451612903225
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = -14.285714285714
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -40
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -28.571428571428
$ws.Range("L28").Value = -16.666666666666
$ws.Range("N29").Value = -96.428571428571
$ws.Range("N30").Value = -96.428571428571
$ws.Range("F31").Value = 2
$ws.Range("I31").Value = 2

# --- Numeric cells converted to text placeholders ("0" or "***.*") ---
# Use a source cell with style 13 (text placeholder style) to copy formats from.
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("C28").Value = "'0"
$ws.Range("F29").Value = "'0"
$ws.Range("F30").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122)

# --- Text placeholder cells converted to real numbers ---
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 3
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 100
$ws.Range("C15").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("H15").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = 100
$ws.Range("C15").Copy() | Out-Null
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 1
$ws.Range("H15").Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("K22").Value = 200
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 1

$excel.CutCopyMode = 0
